$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.026113033294678
$ws.Range("B1").Value = 6.359650611877441
$ws.Range("C1").Value = 7.085855960845947
$ws.Range("D1").Value = 7.611157894134521
$ws.Range("E1").Value = 4.642078399658203
